# Database_Thresholds.xlsx update: refresh the indicator-quantile script /
# details-workbook provenance strings recorded on each threshold row.
#
# - Column W (ScriptLatestRunVersion) on every data row (4-92) is stamped
#   with the new IndicatorQuantiles.R commit id.
# - Column U (QuantileSource) on the handful of rows that were sourced from
#   Database_Thresholds_details.xlsx with the *old* commit id get bumped to
#   the new commit id (rows 29, 41, 62, 63, 78 - identified by ActionNeeded
#   ("U") rows whose QuantileSource points at that workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newScriptVersion = "IndicatorQuantiles.R, Git Commit ID: db49f0f869e1f5a8558dc746458075a467cf2c41"
$newDetailsSource  = "Database_Thresholds_details.xlsx, Git Commit ID: 17b6a0f858dccbb28fc8ab3fe179e7fa731e5996"

$firstDataRow = 4
$lastDataRow  = 92

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $ws.Range("W$row").Value = $newScriptVersion
}

$detailsRows = @(29, 41, 62, 63, 78)
foreach ($row in $detailsRows) {
    $ws.Range("U$row").Value = $newDetailsSource
}
